$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the date number format used by the old date column before
# removing the stale rows.
$dateFormat = $ws.Range("A2").NumberFormat

# Remove old data rows (rows 2-9), keeping the header row (row 1) intact.
$ws.Range("A2:C9").EntireRow.Delete()

# Write the new single data row, re-applying the original date formatting.
$ws.Range("A2").Value = 45426
$ws.Range("A2").NumberFormat = $dateFormat
$ws.Range("B2").Value = "Salary"
$ws.Range("C2").Value = 2000

# Match the author's final selection (top-left cell).
$ws.Range("A1").Select() | Out-Null
